$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly data: rotate the recorded price-report rows (row 3 and row 10 stay put).
# Values below reflect the updated dataset per the diff.

# Row 2
$ws.Range("D2").Value = 45141
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 8500
$ws.Range("L2").Value = 9000
$ws.Range("M2").Value = 8800
$ws.Range("P2").Value = 587

# Row 4
$ws.Range("D4").Value = 45119
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 20000
$ws.Range("L4").Value = 20000
$ws.Range("M4").Value = 20000
$ws.Range("P4").Value = 1333

# Row 5
$ws.Range("D5").Value = 45084
$ws.Range("J5").Value = 90
$ws.Range("K5").Value = 22000
$ws.Range("L5").Value = 23000
$ws.Range("M5").Value = 22556
$ws.Range("P5").Value = 1504

# Row 6
$ws.Range("D6").Value = 44749
$ws.Range("J6").Value = 90
$ws.Range("K6").Value = 17000
$ws.Range("L6").Value = 18000
$ws.Range("M6").Value = 17556
$ws.Range("P6").Value = 1170

# Row 7
$ws.Range("D7").Value = 45133
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 22000
$ws.Range("L7").Value = 22000
$ws.Range("M7").Value = 22000
$ws.Range("P7").Value = 1467

# Row 8
$ws.Range("D8").Value = 44750
$ws.Range("J8").Value = 140
$ws.Range("K8").Value = 19000
$ws.Range("L8").Value = 20000
$ws.Range("M8").Value = 19571
$ws.Range("P8").Value = 1305

# Row 9
$ws.Range("D9").Value = 45091
$ws.Range("J9").Value = 40
$ws.Range("K9").Value = 20000
$ws.Range("L9").Value = 22000
$ws.Range("M9").Value = 21000
$ws.Range("P9").Value = 1400
